$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("F8").Value = "[[Cyanide]]"
